# Dodan opis k projektu
# Shorten the crop-name labels, drop the "Deteljno travne mešanice (1 do 5 let)"
# data row, and re-point the selection, matching the author's edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the whole data row for "Deteljno travne mešanice (1 do 5 let)" - this
# shifts every following row up by one and also drops the now-unused shared
# string automatically on save.
$ws.Rows(15).Delete()

# Replace the long table title / crop-name labels with their shortened forms.
$ws.Range("A1").Value  = "pridelek"
$ws.Range("A2").Value  = "Pšenica "
$ws.Range("A3").Value  = "Ječmen"
$ws.Range("A4").Value  = "Tritikala"
$ws.Range("A5").Value  = "Oves"
$ws.Range("A6").Value  = "Riž"
$ws.Range("A7").Value  = "Koruza "
$ws.Range("A8").Value  = "Silažna"
$ws.Range("A9").Value  = "Krompir"
$ws.Range("A10").Value = "Buče za olje"
$ws.Range("A11").Value = "Repica"
$ws.Range("A12").Value = "Hmelj"
$ws.Range("A13").Value = "Trave"
$ws.Range("A14").Value = "Deteljne"
$ws.Range("A15").Value = "Detelja"
$ws.Range("A16").Value = "Lucerna"
$ws.Range("A17").Value = "Trajni travniki"
$ws.Range("A18").Value = "Zelje"
$ws.Range("A19").Value = "Grozdje"
$ws.Range("A20").Value = "Jabolka "
$ws.Range("A21").Value = "Oljke"
$ws.Range("A22").Value = "Breskve"

# Match the author's final cursor position.
[void]$ws.Range("A25").Select()
